$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 loses the bold/header style applied to it, back to Normal
$ws.Range("A2").Style = "Normal"

# I2 becomes a real number (was text "2")
$ws.Range("I2").Value = 2

# New row 3: N_Ano = 2, rest blank
$ws.Range("A3").Value = 2

# New row 4: N_Ano = 3, rest blank
$ws.Range("A4").Value = 3

Write-Output "done"
